# Update NATMI ligand-receptor pair (Lipc-Lrp1) stats with new TPM-derived
# expression values. Recomputed downstream specificity/edge-weight columns
# for all 9 sending/target cluster combinations accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07453733333333333
$ws.Range("H2").Value = 0.223612
$ws.Range("I2").Value = 0.3946997283496047
$ws.Range("J2").Value = 0.3946997283496047
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 0.2576208012391111
$ws.Range("R2").Value = 2.318587211152
$ws.Range("S2").Value = 0.003884451509833217
$ws.Range("T2").Value = 0.003884451509833216
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07453733333333333
$ws.Range("H3").Value = 0.223612
$ws.Range("I3").Value = 0.3946997283496047
$ws.Range("J3").Value = 0.3946997283496047
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 22.48047081455022
$ws.Range("R3").Value = 202.324237330952
$ws.Range("S3").Value = 0.338964471724824
$ws.Range("T3").Value = 0.3389644717248239
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07453733333333333
$ws.Range("H4").Value = 0.223612
$ws.Range("I4").Value = 0.3946997283496047
$ws.Range("J4").Value = 0.3946997283496047
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 3.438798482820889
$ws.Range("R4").Value = 30.949186345388
$ws.Range("S4").Value = 0.05185080511494752
$ws.Range("T4").Value = 0.05185080511494751
$ws.Range("I5").Value = 0.3614891878200364
$ws.Range("J5").Value = 0.3614891878200364
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 0.2359442571568889
$ws.Range("R5").Value = 2.123498314412
$ws.Range("S5").Value = 0.003557608786019146
$ws.Range("T5").Value = 0.003557608786019145
$ws.Range("I6").Value = 0.3614891878200364
$ws.Range("J6").Value = 0.3614891878200364
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.3104435670528808
$ws.Range("T6").Value = 0.3104435670528807
$ws.Range("I7").Value = 0.3614891878200364
$ws.Range("J7").Value = 0.3614891878200364
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.04748801198113655
$ws.Range("T7").Value = 0.04748801198113654
$ws.Range("I8").Value = 0.2438110838303588
$ws.Range("J8").Value = 0.2438110838303588
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 0.1591356726542222
$ws.Range("R8").Value = 1.432221053888
$ws.Range("S8").Value = 0.002399475511825137
$ws.Range("T8").Value = 0.002399475511825137
$ws.Range("I9").Value = 0.2438110838303588
$ws.Range("J9").Value = 0.2438110838303588
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.2093827010643726
$ws.Range("T9").Value = 0.2093827010643726
$ws.Range("I10").Value = 0.2438110838303588
$ws.Range("J10").Value = 0.2438110838303588
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("S10").Value = 0.0320289072541611
$ws.Range("T10").Value = 0.03202890725416109
